$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add a new "2022" column (S) mirroring the existing year columns ---
# Write the values first, then copy the formatting-only (xlPasteFormats = -4122)
# from the neighbouring cells so the new cells pick up the same styles
# (s="21" for the year header row, s="22" for the computed growth-rate row)
# without disturbing the values we just wrote.

$ws.Range("S4").Value = 2022
$ws.Range("S5").Value = 109.27053140096621

# --- Update the recomputed growth-rate figures for 2020 and 2021 ---
$ws.Range("Q5").Value = 91.892815141492093
$ws.Range("R5").Value = 101.53074848578628

# Copy formatting from R4 (year header) onto the new S4 header cell
$ws.Range("R4").Copy()
$ws.Range("S4").PasteSpecial(-4122)

# Copy formatting from R5 (computed-value style) onto S5 and Q5
$ws.Range("R5").Copy()
$ws.Range("S5").PasteSpecial(-4122)
$ws.Range("Q5").PasteSpecial(-4122)

# --- Update the active cell / selection shown in the worksheet ---
$ws.Range("T5").Select()
